# The post "「私は〜する」" (row 681) was removed from the posts sheet.
# Deleting the entire row shifts every following row (682..857) up by one,
# which matches the target diff (rows 682-857 become 681-856 and the
# sheet dimension shrinks from A1:C857 to A1:C856).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(681).Delete()
